# Add new columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column headers, matching style of existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-31
$data = @(
    @(5, 5),
    @(8, 8),
    @(5, 5),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(5, 7),
    @(8, 8),
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(5, 6),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(8, 9),
    @(5, 6),
    @(5, 5),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
